## Weekly update: insert two new price rows (Sandia / Terminal Hortofrutícola
## Agro Chillán) above the existing daily records, pushing the existing
## rows 96-162 down to 98-164.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 96-97; everything that was on rows 96.. shifts
# down by two (old row 96 -> 98, ... old row 162 -> 164).
$ws.Range("96:97").Insert()

# New row 96: Sandia, Extra
$ws.Range("A96").Value = 7
$ws.Range("B96").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C96").Value = "Ñuble"
$ws.Range("D96").Value = 44596
$ws.Range("E96").Value = 16
$ws.Range("F96").Value = 100112028
$ws.Range("G96").Value = "Sandia"
$ws.Range("H96").Value = "Sin especificar"
$ws.Range("I96").Value = "Extra"
$ws.Range("J96").Value = 500
$ws.Range("K96").Value = 2500
$ws.Range("L96").Value = 2500
$ws.Range("M96").Value = 2500
$ws.Range("N96").Value = "`$/unidad"
$ws.Range("O96").Value = "Región de O'Higgins"
$ws.Range("P96").Value = 2500
$ws.Range("Q96").Value = 1
$ws.Range("R96").Value = "Hortaliza"

# New row 97: Sandia, Primera
$ws.Range("A97").Value = 7
$ws.Range("B97").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C97").Value = "Ñuble"
$ws.Range("D97").Value = 44596
$ws.Range("E97").Value = 16
$ws.Range("F97").Value = 100112028
$ws.Range("G97").Value = "Sandia"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 600
$ws.Range("K97").Value = 1800
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = 1900
$ws.Range("N97").Value = "`$/unidad"
$ws.Range("O97").Value = "Región de O'Higgins"
$ws.Range("P97").Value = 1900
$ws.Range("Q97").Value = 1
$ws.Range("R97").Value = "Hortaliza"
